# Updated cryptos list on Mon Oct 16 05:54:05 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row.
# Price/volume values are free-form text (e.g. "27.540.69", padded "  +1.97%  ")
# rather than numbers, so plain numeric-looking prices are forced to stay text
# (NumberFormat "@" while writing, then restored to the "Normal" style so no
# lasting formatting change is left behind).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.540.69"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "1.570.49"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.55"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.250"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("D12").Value = "1.798.43"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "1.588.29"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "27.511.34"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.991"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.459.28"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("E35").Value = "  +4.04%  "
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.542"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.817"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.966"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.35%  "
$ws.Range("E46").Value = "  +3.44%  "
$ws.Range("D47").Value = "1.711.31"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").Value = "0.0₇0957"
$ws.Range("E50").Value = "  -7.58%  "
$ws.Range("E51").Value = "  -2.05%  "
